$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18: dct:hasVersion
$ws.Range("A18").Value = "dct:hasVersion^^xsd:string"
$ws.Range("B18").Value = "0.1.0"

# Row 19: dct:created
$ws.Range("A19").Value = "dct:created^^xsd:datetime"
$ws.Range("B19").Value = "2022-06-01T00:00:00+00:00"

# Row 20: dct:modified
$ws.Range("A20").Value = "dct:modified^^xsd:datetime"
$ws.Range("B20").Value = "2022-06-11T18:35:01+00:00"

$wb.Save()
